# Apply changes described by the diff:
# 1. Rename the worksheet from "Hoja1" to "REFRESCOS PREFERIDOS"
# 2. Insert a new header row "refresco" at A1, shifting the existing
#    100 values (A1:A100) down to A2:A101
# 3. Update the selection to C4 (no frozen/scrolled topLeftCell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "REFRESCOS PREFERIDOS"

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Set the new header value
$ws.Range("A1").Value = "refresco"

# Update the selection/view
$ws.Range("C4").Select()
